# Add 9 new "Metal" katana pre-order rows to the Armory section.
# This inserts 9 rows above the current row 20 (shifting all rows from
# 20 onward down by 9) and fills them with the new product data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows at row 20 (existing data shifts down to rows 29..85).
$ws.Rows("20:28").Insert()

# Column A - product names (filled first, row by row, so the new shared
# strings are appended to the string table in this order).
$ws.Range("A20").Value = "Shusui Katana [Metal]"
$ws.Range("A21").Value = "Wado Ichimonji Katana [Metal]"
$ws.Range("A22").Value = "Enma Katana [Metal]"
$ws.Range("A23").Value = "Enma Black Katana [Metal]"
$ws.Range("A24").Value = "Ame No Habikiri Katana [Metal]"
$ws.Range("A25").Value = "Tanjiro's Old Nichrin [Metal]"
$ws.Range("A26").Value = "Tanjiro's New Nichrin [Metal]"
$ws.Range("A27").Value = "Inosuke's Nichrin [Metal]"
$ws.Range("A28").Value = "Sanemi's Nichrin [Metal]"

# Column B - price.
$ws.Range("B20").Value = 5700
$ws.Range("B21").Value = 5700
$ws.Range("B22").Value = 5700
$ws.Range("B23").Value = 5700
$ws.Range("B24").Value = 5700
$ws.Range("B25").Value = 5700
$ws.Range("B26").Value = 5700
$ws.Range("B27").Value = 5700
$ws.Range("B28").Value = 5700

# Column C - category (reuses the existing "Armory" shared string).
$ws.Range("C20").Value = "Armory"
$ws.Range("C21").Value = "Armory"
$ws.Range("C22").Value = "Armory"
$ws.Range("C23").Value = "Armory"
$ws.Range("C24").Value = "Armory"
$ws.Range("C25").Value = "Armory"
$ws.Range("C26").Value = "Armory"
$ws.Range("C27").Value = "Armory"
$ws.Range("C28").Value = "Armory"

# Column D - image file names. Row 24's image is written before row 23's
# so the new shared strings land in the exact index order of the target
# workbook (...,172,174,173,175,...).
$ws.Range("D20").Value = "shusui metal.jpg"
$ws.Range("D21").Value = "wado metal.jpg"
$ws.Range("D22").Value = "enma metal.jpg"
$ws.Range("D24").Value = "ame metal.jpg"
$ws.Range("D23").Value = "enma bl metal.jpg"
$ws.Range("D25").Value = "tanjiro old metal.jpg"
$ws.Range("D26").Value = "tanjiro new metal.jpg"
$ws.Range("D27").Value = "ino metal.jpg"
$ws.Range("D28").Value = "sanemi metal.jpg"

# Update the view state to match: scroll so row 13 is at the top and
# select D28.
$ws.Range("D28").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
